$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Execute" column (C) to "Y" for rows 2-12 and refresh the
# "Date" column (B) with the new run's timestamps.
$ws.Range("C2").Value = "Y"
$ws.Range("B2").Value = "Tue Oct 25 22:43:44 EDT 2022"

$ws.Range("C3").Value = "Y"
$ws.Range("B3").Value = "Tue Oct 25 22:44:24 EDT 2022"

$ws.Range("C4").Value = "Y"
$ws.Range("B4").Value = "Tue Oct 25 22:45:00 EDT 2022"

$ws.Range("C5").Value = "Y"
$ws.Range("B5").Value = "Tue Oct 25 22:46:10 EDT 2022"

$ws.Range("C6").Value = "Y"
$ws.Range("B6").Value = "Tue Oct 25 22:46:51 EDT 2022"

$ws.Range("C7").Value = "Y"
$ws.Range("B7").Value = "Tue Oct 25 22:47:28 EDT 2022"

$ws.Range("C8").Value = "Y"
$ws.Range("B8").Value = "Tue Oct 25 22:48:05 EDT 2022"

$ws.Range("C9").Value = "Y"
$ws.Range("B9").Value = "Tue Oct 25 22:48:42 EDT 2022"

$ws.Range("C10").Value = "Y"
$ws.Range("B10").Value = "Tue Oct 25 22:49:19 EDT 2022"

$ws.Range("C11").Value = "Y"
$ws.Range("B11").Value = "Tue Oct 25 22:49:56 EDT 2022"

$ws.Range("C12").Value = "Y"
$ws.Range("B12").Value = "Tue Oct 25 22:50:34 EDT 2022"

# Selection moved from E18 to C2:C12.
$ws.Range("C2:C12").Select()
